$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 906.2
$ws.Range("I4").Value = 920.8889
$ws.Range("K4").Value = 920.8889
$ws.Range("M4").Value = -806.8889
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H107").Value = 1238.174
$ws.Range("I107").Value = 1244.3334
$ws.Range("K107").Value = 1244.3334
$ws.Range("M107").Value = 675.6666
$ws.Range("H117").Value = 0
$ws.Range("J117").Value = 0
$ws.Range("L117").Value = 0
$ws.Range("N117").ClearContents()
$ws.Range("H138").Value = 3667.68
$ws.Range("J138").Value = 5160.6
$ws.Range("L138").Value = 15481.8
$ws.Range("N138").Value = -25761.8
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12504.091
$ws.Range("I2").Value = 14677.625
$ws.Range("K2").Value = 14677.625
$ws.Range("M2").Value = -14564.625
$ws.Range("H5").Value = 57.5
$ws.Range("I5").Value = 57.5
$ws.Range("K5").Value = 57.5
$ws.Range("M5").Value = 54.5
$ws.Range("H45").Value = 1122.4
$ws.Range("J45").Value = 1500
$ws.Range("L45").Value = 1500
$ws.Range("N45").Value = -2254
$ws.Range("H110").Value = 3224.3914
$ws.Range("I110").Value = 2650.0588
$ws.Range("K110").Value = 2650.0588
$ws.Range("M110").Value = -605.0587999999998
$ws.Range("H116").Value = 12504.091
$ws.Range("I116").Value = 14677.625
$ws.Range("K116").Value = 14677.625
$ws.Range("M116").Value = -12383.625
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 12504.091
$ws.Range("I3").Value = 14677.625
$ws.Range("K3").Value = 14677.625
$ws.Range("M3").Value = -14563.625
$ws.Range("H4").Value = 57.5
$ws.Range("I4").Value = 57.5
$ws.Range("K4").Value = 57.5
$ws.Range("M4").Value = 57.5
$ws.Range("H53").Value = 517500
$ws.Range("J53").Value = 517500
$ws.Range("L53").Value = 517500
$ws.Range("N53").Value = -518648
$ws.Range("H86").Value = 5699.25
$ws.Range("J86").Value = 5699.25
$ws.Range("L86").Value = 5699.25
$ws.Range("N86").Value = -7945.25
$ws.Range("H89").Value = 5699.25
$ws.Range("J89").Value = 5699.25
$ws.Range("L89").Value = 28496.25
$ws.Range("N89").Value = -39728.25
$ws.Range("H107").Value = 9014.588
$ws.Range("I107").Value = 9014.588
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 9014.588
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -7094.588
$ws.Range("N107").ClearContents()
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 3242
$ws.Range("I58").Value = 3125.875
$ws.Range("J58").Value = 3706.5
$ws.Range("K58").Value = 3125.875
$ws.Range("L58").Value = 3706.5
$ws.Range("M58").Value = -2922.875
$ws.Range("N58").Value = -4112.5
$ws.Range("H86").Value = 10657.25
$ws.Range("I86").Value = 5849.5
$ws.Range("K86").Value = 5849.5
$ws.Range("M86").Value = -4726.5
$ws.Range("H89").Value = 10657.25
$ws.Range("I89").Value = 5849.5
$ws.Range("K89").Value = 29247.5
$ws.Range("M89").Value = -23631.5
$ws.Range("H94").Value = 1595
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1595
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1595
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -2497
$ws.Range("H107").Value = 3183.4546
$ws.Range("I107").Value = 682.6667
$ws.Range("K107").Value = 682.6667
$ws.Range("M107").Value = 1237.3333
$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()
$ws.Range("H136").Value = 3242
$ws.Range("I136").Value = 3125.875
$ws.Range("J136").Value = 3706.5
$ws.Range("K136").Value = 9377.625
$ws.Range("L136").Value = 11119.5
$ws.Range("M136").Value = -6827.625
$ws.Range("N136").Value = -16219.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 2259
$ws.Range("I81").Value = 149
$ws.Range("J81").Value = 3665.6667
$ws.Range("K81").Value = 447
$ws.Range("L81").Value = 10997.0001
$ws.Range("M81").Value = 676
$ws.Range("N81").Value = -13243.0001
$ws.Range("H84").Value = 2259
$ws.Range("I84").Value = 149
$ws.Range("J84").Value = 3665.6667
$ws.Range("K84").Value = 1341
$ws.Range("L84").Value = 32991.0003
$ws.Range("M84").Value = 4275
$ws.Range("N84").Value = -44223.0003
$ws.Range("H88").Value = 3902
$ws.Range("I88").Value = 3869.6667
$ws.Range("K88").Value = 11609.0001
$ws.Range("M88").Value = -11181.0001
$ws.Range("H91").Value = 3902
$ws.Range("I91").Value = 3869.6667
$ws.Range("K91").Value = 11609.0001
$ws.Range("M91").Value = -10127.0001
$ws.Range("H107").Value = 1482.0741
$ws.Range("I107").Value = 225
$ws.Range("J107").Value = 1530.4231
$ws.Range("K107").Value = 675
$ws.Range("L107").Value = 4591.2693
$ws.Range("M107").Value = 1245
$ws.Range("N107").Value = -8431.2693
$ws.Range("H118").Value = 1899.875
$ws.Range("I118").Value = 1899.875
$ws.Range("K118").Value = 5699.625
$ws.Range("M118").Value = -4456.625
$ws.Range("H125").Value = 7400
$ws.Range("I125").Value = 7400
$ws.Range("K125").Value = 22200
$ws.Range("M125").Value = -17280
$ws.Range("H140").Value = 2822.6924
$ws.Range("I140").Value = 1485.3158
$ws.Range("K140").Value = 4455.9474
$ws.Range("M140").Value = 724.0526
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2795.182
$ws.Range("J80").Value = 3056.8572
$ws.Range("L80").Value = 3056.8572
$ws.Range("N80").Value = -5052.8572
$ws.Range("H83").Value = 2795.182
$ws.Range("J83").Value = 3056.8572
$ws.Range("L83").Value = 15284.286
$ws.Range("N83").Value = -25268.286
$ws.Range("H97").Value = 1142.5238
$ws.Range("I97").Value = 1131.4706
$ws.Range("K97").Value = 1131.4706
$ws.Range("M97").Value = -635.4706000000001
$ws.Range("H132").Value = 3840.5715
$ws.Range("I132").Value = 3647.8333
$ws.Range("J132").Value = 4997
$ws.Range("K132").Value = 10943.4999
$ws.Range("L132").Value = 14991
$ws.Range("M132").Value = -8413.499899999999
$ws.Range("N132").Value = -20051
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 80010
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H7").Value = 5000
$ws.Range("I7").Value = 5000
$ws.Range("K7").Value = 5000
$ws.Range("M7").Value = -4888
$ws.Range("H28").Value = 80010
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()
$ws.Range("H37").Value = 80010
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()
$ws.Range("H40").Value = 5999.125
$ws.Range("I40").Value = 5065.8335
$ws.Range("K40").Value = 5065.8335
$ws.Range("M40").Value = -4929.8335
$ws.Range("H126").Value = 5000
$ws.Range("I126").Value = 5000
$ws.Range("K126").Value = 15000
$ws.Range("M126").Value = -12530
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H24").Value = 36600.2
$ws.Range("J24").Value = 40501.75
$ws.Range("L24").Value = 40501.75
$ws.Range("N24").Value = -40961.75
$ws.Range("H107").Value = 1728.5
$ws.Range("I107").Value = 1983
$ws.Range("K107").Value = 5949
$ws.Range("M107").Value = -4029
$ws.Range("H122").Value = 2448.2666
$ws.Range("I122").Value = 2320
$ws.Range("K122").Value = 6960
$ws.Range("M122").Value = -4510
$ws.Range("H126").Value = 3150.625
$ws.Range("I126").Value = 2529.3572
$ws.Range("K126").Value = 7588.071599999999
$ws.Range("M126").Value = -5118.071599999999
$ws.Range("H136").Value = 1628.1515
$ws.Range("I136").Value = 1538.4062
$ws.Range("K136").Value = 4615.2186
$ws.Range("M136").Value = -2065.2186
